# Apply updated "想去人数" (F) / "最低票价" (G) figures across the four sheets
# of the 上海-漫展信息 workbook, matching the refreshed data pull.

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsAll = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibitions)
$wsExpo.Range("F6").Value = 1326
$wsExpo.Range("F12").Value = 1817
$wsExpo.Range("F18").Value = 527
$wsExpo.Range("F22").Value = 17
$wsExpo.Range("F23").Value = 2414
$wsExpo.Range("F24").Value = 458
$wsExpo.Range("F26").Value = 1045
$wsExpo.Range("F27").Value = 4614
$wsExpo.Range("F34").Value = 1004

# 演出 (Performances)
$wsShow.Range("F2").Value = 34
$wsShow.Range("F8").Value = 174
$wsShow.Range("F10").Value = 3
$wsShow.Range("F39").Value = 47
$wsShow.Range("F42").Value = 103
$wsShow.Range("G42").Value = 180

# 本地生活 (Local life)
$wsLocal.Range("F9").Value = 3126
$wsLocal.Range("F10").Value = 629
$wsLocal.Range("F11").Value = 899
$wsLocal.Range("F14").Value = 63
$wsLocal.Range("F15").Value = 17
$wsLocal.Range("F16").Value = 319

# 全部类型 (All types - aggregated view)
$wsAll.Range("F6").Value = 3126
$wsAll.Range("F7").Value = 629
$wsAll.Range("F8").Value = 899
$wsAll.Range("F12").Value = 63
$wsAll.Range("F13").Value = 63
$wsAll.Range("F14").Value = 1326
$wsAll.Range("F17").Value = 17
$wsAll.Range("F18").Value = 1817
$wsAll.Range("F23").Value = 527
$wsAll.Range("F29").Value = 17
$wsAll.Range("F31").Value = 2414
$wsAll.Range("F32").Value = 458
$wsAll.Range("F35").Value = 1045
$wsAll.Range("F37").Value = 319
$wsAll.Range("F44").Value = 47
$wsAll.Range("F47").Value = 103
$wsAll.Range("G47").Value = 180
$wsAll.Range("F51").Value = 1004
